$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    # Force text storage so numeric-looking strings (e.g. "1.00") are not
    # silently coerced into numbers, then restore the original (General)
    # number format so no stray style is left behind on the cell.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

$updates = @(
    @{Addr='D2'; Val='98.473.74'},
    @{Addr='E2'; Val='  +0.88%  '},
    @{Addr='D3'; Val='3.335.84'},
    @{Addr='E3'; Val='  +6.39%  '},
    @{Addr='E4'; Val='  +0.00%  '},
    @{Addr='D5'; Val='258.43'},
    @{Addr='E5'; Val='  +6.84%  '},
    @{Addr='D6'; Val='626.71'},
    @{Addr='E6'; Val='  +2.71%  '},
    @{Addr='D7'; Val='1.40'},
    @{Addr='E7'; Val='  +24.63%  '},
    @{Addr='D8'; Val='0.388'},
    @{Addr='E8'; Val='  +1.24%  '},
    @{Addr='D9'; Val='1.00'},
    @{Addr='E9'; Val='  -0.03%  '},
    @{Addr='D10'; Val='0.864'},
    @{Addr='E10'; Val='  +10.57%  '},
    @{Addr='D11'; Val='3.334.69'},
    @{Addr='E11'; Val='  +6.39%  '},
    @{Addr='D12'; Val='0.199'},
    @{Addr='B13'; Val='Avalanche'},
    @{Addr='C13'; Val='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'},
    @{Addr='D13'; Val='37.06'},
    @{Addr='E13'; Val='  +9.94%  '},
    @{Addr='B14'; Val='WrappedBTC'},
    @{Addr='C14'; Val='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'},
    @{Addr='D14'; Val='98.227.12'},
    @{Addr='E14'; Val='  +0.96%  '},
    @{Addr='D15'; Val='0.0000248'},
    @{Addr='E15'; Val='  +3.96%  '},
    @{Addr='D16'; Val='3.979.54'},
    @{Addr='E16'; Val='  +7.19%  '},
    @{Addr='D17'; Val='5.50'},
    @{Addr='E17'; Val='  +1.01%  '},
    @{Addr='D18'; Val='3.338.76'},
    @{Addr='E18'; Val='  +6.68%  '},
    @{Addr='D19'; Val='3.54'},
    @{Addr='E19'; Val='  +3.55%  '},
    @{Addr='D20'; Val='15.14'},
    @{Addr='E20'; Val='  +5.31%  '},
    @{Addr='D21'; Val='488.74'},
    @{Addr='E21'; Val='  -5.56%  '},
    @{Addr='D22'; Val='6.05'},
    @{Addr='E22'; Val='  +7.67%  '},
    @{Addr='D23'; Val='0.0000211'},
    @{Addr='E23'; Val='  +10.57%  '},
    @{Addr='D24'; Val='9.40'},
    @{Addr='E24'; Val='  +8.05%  '},
    @{Addr='D25'; Val='5.63'},
    @{Addr='E25'; Val='  +3.51%  '},
    @{Addr='D26'; Val='88.68'},
    @{Addr='E26'; Val='  +0.46%  '},
    @{Addr='D27'; Val='11.89'},
    @{Addr='E27'; Val='  +3.26%  '},
    @{Addr='D28'; Val='3.512.23'},
    @{Addr='E28'; Val='  +6.66%  '},
    @{Addr='D29'; Val='0.284'},
    @{Addr='E29'; Val='  +17.95%  '},
    @{Addr='E30'; Val='  -0.08%  '},
    @{Addr='E31'; Val='  +10.17%  '},
    @{Addr='E32'; Val='  +12.25%  '},
    @{Addr='D33'; Val='0.997'},
    @{Addr='E33'; Val='  +0.06%  '},
    @{Addr='E34'; Val='  +8.14%  '},
    @{Addr='E35'; Val='  +5.12%  '},
    @{Addr='E36'; Val='  -0.38%  '},
    @{Addr='D37'; Val='7.28'},
    @{Addr='E37'; Val='  +0.75%  '},
    @{Addr='E38'; Val='  +4.15%  '},
    @{Addr='D39'; Val='496.67'},
    @{Addr='E39'; Val='  +6.42%  '},
    @{Addr='D40'; Val='0.460'},
    @{Addr='E40'; Val='  +6.19%  '},
    @{Addr='E41'; Val='  +2.31%  '},
    @{Addr='E42'; Val='  +4.73%  '},
    @{Addr='E43'; Val='  +3.88%  '},
    @{Addr='D44'; Val='3.30'},
    @{Addr='E44'; Val='  +6.60%  '},
    @{Addr='E45'; Val='  +0.00%  '},
    @{Addr='D46'; Val='0.781'},
    @{Addr='E46'; Val='  +12.52%  '},
    @{Addr='D47'; Val='159.45'},
    @{Addr='E47'; Val='  -1.00%  '},
    @{Addr='D48'; Val='1.93'},
    @{Addr='E48'; Val='  +1.85%  '},
    @{Addr='D49'; Val='0.845'},
    @{Addr='E49'; Val='  +9.77%  '},
    @{Addr='D50'; Val='4.60'},
    @{Addr='E50'; Val='  +3.33%  '},
    @{Addr='D51'; Val='45.56'},
    @{Addr='E51'; Val='  +3.33%  '},
)

foreach ($u in $updates) {
    Set-TextCell $ws $u.Addr $u.Val
}
